# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3180
$wsExpo.Range("F4").Value = 989
$wsExpo.Range("F5").Value = 303

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3180
$wsAll.Range("F4").Value = 989
$wsAll.Range("F6").Value = 303
